$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Remove all existing hyperlinks up front -- we will rebuild them after the
# row data has been shifted, since the hyperlink objects do not follow the
# cell-content shift performed below.
$ws.Hyperlinks.Delete()

# Shift the existing data rows (2-16) down by one row (to 3-17) to make room
# for the new, most-recent price entry at row 2. Processing bottom-to-top so
# a source row is never overwritten before it has been copied. Range.Copy
# carries the original cell formatting (style 3 / 4) along with it, so the
# shifted rows keep their look untouched.
for ($r = 16; $r -ge 2; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":F" + $srcRow)
    $dst = $ws.Range("A" + $dstRow)
    $src.Copy($dst)
}

# Populate the new row 2 with the latest circular entry. The date column is
# forced to text format before the value is entered so the dd-mm-yyyy string
# is not auto-coerced into a date serial value; the original look (style 3,
# General number format) is then restored via a formats-only paste from a
# neighbouring cell that already carries that style, avoiding creation of a
# divergent style record.
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value2 = 16
$ws.Range("B2").Value2 = "ALUMINIUM INGOT"
$ws.Range("C2").Value2 = "IE07"
$ws.Range("D2").Value2 = 297.15
$ws.Range("E2").Value2 = "01-11-2025"
$ws.Range("F2").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Rebuild the hyperlinks for column F, rows 2-17, in order. Adding a
# hyperlink re-styles the cell with the built-in "Hyperlink" look, so the
# original style (3) is restored afterwards the same way as above.
$links = @(
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
)

for ($i = 0; $i -lt $links.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("F" + $row)
    $ws.Hyperlinks.Add($cell, $links[$i]) | Out-Null

    $ws.Range("A" + $row).Copy()
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

$ws.Range("A1").Select()
